$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 61.84465033333333
$ws.Range("H2").Value = 185.533951
$ws.Range("I2").Value = 0.03153184209101587
$ws.Range("J2").Value = 0.03153184209101587
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 22.34478233333333
$ws.Range("N2").Value = 67.034347
$ws.Range("O2").Value = 0.03165884810812076
$ws.Range("P2").Value = 0.03165884810812076
$ws.Range("Q2").Value = 1381.905250179444
$ws.Range("R2").Value = 12437.147251615
$ws.Range("S2").Value = 0.0009982617993287205
$ws.Range("T2").Value = 0.0009982617993287203
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 61.84465033333333
$ws.Range("H3").Value = 185.533951
$ws.Range("I3").Value = 0.03153184209101587
$ws.Range("J3").Value = 0.03153184209101587
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 92.44713066666667
$ws.Range("N3").Value = 277.341392
$ws.Range("O3").Value = 0.1309822411400946
$ws.Range("P3").Value = 0.1309822411400946
$ws.Range("Q3").Value = 5717.360470399977
$ws.Range("R3").Value = 51456.24423359979
$ws.Range("S3").Value = 0.004130111344356827
$ws.Range("T3").Value = 0.004130111344356825
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 61.84465033333333
$ws.Range("H4").Value = 185.533951
$ws.Range("I4").Value = 0.03153184209101587
$ws.Range("J4").Value = 0.03153184209101587
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 243.96462
$ws.Range("N4").Value = 731.89386
$ws.Range("O4").Value = 0.3456573768818275
$ws.Range("P4").Value = 0.3456573768818275
$ws.Range("Q4").Value = 15087.90661760454
$ws.Range("R4").Value = 135791.1595584409
$ws.Range("S4").Value = 0.01089921382543254
$ws.Range("T4").Value = 0.01089921382543254
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 61.84465033333333
$ws.Range("H5").Value = 185.533951
$ws.Range("I5").Value = 0.03153184209101587
$ws.Range("J5").Value = 0.03153184209101587
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 281.5837096666667
$ws.Range("N5").Value = 844.751129
$ws.Range("O5").Value = 0.3989573834764815
$ws.Range("P5").Value = 0.3989573834764815
$ws.Range("Q5").Value = 17414.44606389785
$ws.Range("R5").Value = 156730.0145750807
$ws.Range("S5").Value = 0.01257986121682528
$ws.Range("T5").Value = 0.01257986121682528
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 61.84465033333333
$ws.Range("H6").Value = 185.533951
$ws.Range("I6").Value = 0.03153184209101587
$ws.Range("J6").Value = 0.03153184209101587
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 65.45872566666667
$ws.Range("N6").Value = 196.376177
$ws.Range("O6").Value = 0.09274415039347572
$ws.Range("P6").Value = 0.09274415039347571
$ws.Range("Q6").Value = 4048.272000120592
$ws.Range("R6").Value = 36434.44800108533
$ws.Range("S6").Value = 0.002924393905072504
$ws.Range("T6").Value = 0.002924393905072503
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 1361.379069
$ws.Range("H7").Value = 4084.137207
$ws.Range("I7").Value = 0.6941067594101231
$ws.Range("J7").Value = 0.6941067594101232
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 22.34478233333333
$ws.Range("N7").Value = 67.034347
$ws.Range("O7").Value = 0.03165884810812076
$ws.Range("P7").Value = 0.03165884810812076
$ws.Range("Q7").Value = 30419.71896996097
$ws.Range("R7").Value = 273777.4707296488
$ws.Range("S7").Value = 0.02197462046698501
$ws.Range("T7").Value = 0.02197462046698501
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 1361.379069
$ws.Range("H8").Value = 4084.137207
$ws.Range("I8").Value = 0.6941067594101231
$ws.Range("J8").Value = 0.6941067594101232
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 92.44713066666667
$ws.Range("N8").Value = 277.341392
$ws.Range("O8").Value = 0.1309822411400946
$ws.Range("P8").Value = 0.1309822411400946
$ws.Range("Q8").Value = 125855.588678708
$ws.Range("R8").Value = 1132700.298108372
$ws.Range("S8").Value = 0.09091565893802639
$ws.Range("T8").Value = 0.09091565893802639
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 1361.379069
$ws.Range("H9").Value = 4084.137207
$ws.Range("I9").Value = 0.6941067594101231
$ws.Range("J9").Value = 0.6941067594101232
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 243.96462
$ws.Range("N9").Value = 731.89386
$ws.Range("O9").Value = 0.3456573768818275
$ws.Range("P9").Value = 0.3456573768818275
$ws.Range("Q9").Value = 332128.3272445388
$ws.Range("R9").Value = 2989154.945200849
$ws.Range("S9").Value = 0.2399231217336489
$ws.Range("T9").Value = 0.2399231217336489
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 1361.379069
$ws.Range("H10").Value = 4084.137207
$ws.Range("I10").Value = 0.6941067594101231
$ws.Range("J10").Value = 0.6941067594101232
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 281.5837096666667
$ws.Range("N10").Value = 844.751129
$ws.Range("O10").Value = 0.3989573834764815
$ws.Range("P10").Value = 0.3989573834764815
$ws.Range("Q10").Value = 383342.1685115729
$ws.Range("R10").Value = 3450079.516604157
$ws.Range("S10").Value = 0.2769190165876024
$ws.Range("T10").Value = 0.2769190165876024
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 1361.379069
$ws.Range("H11").Value = 4084.137207
$ws.Range("I11").Value = 0.6941067594101231
$ws.Range("J11").Value = 0.6941067594101232
$ws.Range("K11").Value = 3.0
$ws.Range("M11").Value = 65.45872566666667
$ws.Range("N11").Value = 196.376177
$ws.Range("O11").Value = 0.09274415039347572
$ws.Range("P11").Value = 0.09274415039347571
$ws.Range("Q11").Value = 89114.13900601307
$ws.Range("R11").Value = 802027.2510541177
$ws.Range("S11").Value = 0.06437434168386053
$ws.Range("T11").Value = 0.06437434168386053
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 251.007014
$ws.Range("H12").Value = 753.021042
$ws.Range("I12").Value = 0.127977334927537
$ws.Range("J12").Value = 0.1279773349275369
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 22.34478233333333
$ws.Range("N12").Value = 67.034347
$ws.Range("O12").Value = 0.03165884810812076
$ws.Range("P12").Value = 0.03165884810812076
$ws.Range("Q12").Value = 5608.697091969952
$ws.Range("R12").Value = 50478.27382772957
$ws.Range("S12").Value = 0.00405161500775299
$ws.Range("T12").Value = 0.004051615007752989
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 251.007014
$ws.Range("H13").Value = 753.021042
$ws.Range("I13").Value = 0.127977334927537
$ws.Range("J13").Value = 0.1279773349275369
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 92.44713066666667
$ws.Range("N13").Value = 277.341392
$ws.Range("O13").Value = 0.1309822411400946
$ws.Range("P13").Value = 0.1309822411400946
$ws.Range("Q13").Value = 23204.87822150783
$ws.Range("R13").Value = 208843.9039935704
$ws.Range("S13").Value = 0.0167627581439453
$ws.Range("T13").Value = 0.01676275814394529
$ws.Range("E14").Value = 3.0
$ws.Range("G14").Value = 251.007014
$ws.Range("H14").Value = 753.021042
$ws.Range("I14").Value = 0.127977334927537
$ws.Range("J14").Value = 0.1279773349275369
$ws.Range("K14").Value = 3.0
$ws.Range("M14").Value = 243.96462
$ws.Range("N14").Value = 731.89386
$ws.Range("O14").Value = 0.3456573768818275
$ws.Range("P14").Value = 0.3456573768818275
$ws.Range("Q14").Value = 61236.83078784468
$ws.Range("R14").Value = 551131.4770906022
$ws.Range("S14").Value = 0.04423630989137951
$ws.Range("T14").Value = 0.04423630989137949
$ws.Range("E15").Value = 3.0
$ws.Range("G15").Value = 251.007014
$ws.Range("H15").Value = 753.021042
$ws.Range("I15").Value = 0.127977334927537
$ws.Range("J15").Value = 0.1279773349275369
$ws.Range("K15").Value = 3.0
$ws.Range("M15").Value = 281.5837096666667
$ws.Range("N15").Value = 844.751129
$ws.Range("O15").Value = 0.3989573834764815
$ws.Range("P15").Value = 0.3989573834764815
$ws.Range("Q15").Value = 70679.48615447294
$ws.Range("R15").Value = 636115.3753902564
$ws.Range("S15").Value = 0.05105750268698347
$ws.Range("T15").Value = 0.05105750268698346
$ws.Range("E16").Value = 3.0
$ws.Range("G16").Value = 251.007014
$ws.Range("H16").Value = 753.021042
$ws.Range("I16").Value = 0.127977334927537
$ws.Range("J16").Value = 0.1279773349275369
$ws.Range("K16").Value = 3.0
$ws.Range("M16").Value = 65.45872566666667
$ws.Range("N16").Value = 196.376177
$ws.Range("O16").Value = 0.09274415039347572
$ws.Range("P16").Value = 0.09274415039347571
$ws.Range("Q16").Value = 16430.59926983516
$ws.Range("R16").Value = 147875.3934285164
$ws.Range("S16").Value = 0.0118691491974757
$ws.Range("T16").Value = 0.0118691491974757
$ws.Range("E17").Value = 3.0
$ws.Range("G17").Value = 260.0315303333334
$ws.Range("H17").Value = 780.094591
$ws.Range("I17").Value = 0.1325785352324417
$ws.Range("J17").Value = 0.1325785352324417
$ws.Range("K17").Value = 3.0
$ws.Range("M17").Value = 22.34478233333333
$ws.Range("N17").Value = 67.034347
$ws.Range("O17").Value = 0.03165884810812076
$ws.Range("P17").Value = 0.03165884810812076
$ws.Range("Q17").Value = 5810.347945101897
$ws.Range("R17").Value = 52293.13150591707
$ws.Range("S17").Value = 0.004197283709321009
$ws.Range("T17").Value = 0.004197283709321009
$ws.Range("E18").Value = 3.0
$ws.Range("G18").Value = 260.0315303333334
$ws.Range("H18").Value = 780.094591
$ws.Range("I18").Value = 0.1325785352324417
$ws.Range("J18").Value = 0.1325785352324417
$ws.Range("K18").Value = 3.0
$ws.Range("M18").Value = 92.44713066666667
$ws.Range("N18").Value = 277.341392
$ws.Range("O18").Value = 0.1309822411400946
$ws.Range("P18").Value = 0.1309822411400946
$ws.Range("Q18").Value = 24039.16886217897
$ws.Range("R18").Value = 216352.5197596107
$ws.Range("S18").Value = 0.01736543367181621
$ws.Range("T18").Value = 0.01736543367181621
$ws.Range("E19").Value = 3.0
$ws.Range("G19").Value = 260.0315303333334
$ws.Range("H19").Value = 780.094591
$ws.Range("I19").Value = 0.1325785352324417
$ws.Range("J19").Value = 0.1325785352324417
$ws.Range("K19").Value = 3.0
$ws.Range("M19").Value = 243.96462
$ws.Range("N19").Value = 731.89386
$ws.Range("O19").Value = 0.3456573768818275
$ws.Range("P19").Value = 0.3456573768818275
$ws.Range("Q19").Value = 63438.49348579015
$ws.Range("R19").Value = 570946.4413721113
$ws.Range("S19").Value = 0.04582674871928074
$ws.Range("T19").Value = 0.04582674871928074
$ws.Range("E20").Value = 3.0
$ws.Range("G20").Value = 260.0315303333334
$ws.Range("H20").Value = 780.094591
$ws.Range("I20").Value = 0.1325785352324417
$ws.Range("J20").Value = 0.1325785352324417
$ws.Range("K20").Value = 3.0
$ws.Range("M20").Value = 281.5837096666667
$ws.Range("N20").Value = 844.751129
$ws.Range("O20").Value = 0.3989573834764815
$ws.Range("P20").Value = 0.3989573834764815
$ws.Range("Q20").Value = 73220.64294156036
$ws.Range("R20").Value = 658985.7864740433
$ws.Range("S20").Value = 0.05289318552147946
$ws.Range("T20").Value = 0.05289318552147945
$ws.Range("E21").Value = 3.0
$ws.Range("G21").Value = 260.0315303333334
$ws.Range("H21").Value = 780.094591
$ws.Range("I21").Value = 0.1325785352324417
$ws.Range("J21").Value = 0.1325785352324417
$ws.Range("K21").Value = 3.0
$ws.Range("M21").Value = 65.45872566666667
$ws.Range("N21").Value = 196.376177
$ws.Range("O21").Value = 0.09274415039347572
$ws.Range("P21").Value = 0.09274415039347571
$ws.Range("Q21").Value = 17021.33260877318
$ws.Range("R21").Value = 153191.9934789586
$ws.Range("S21").Value = 0.01229588361054429
$ws.Range("T21").Value = 0.01229588361054429
$ws.Range("E22").Value = 3.0
$ws.Range("G22").Value = 27.07732933333333
$ws.Range("H22").Value = 81.231988
$ws.Range("I22").Value = 0.01380552833888228
$ws.Range("J22").Value = 0.01380552833888228
$ws.Range("K22").Value = 3.0
$ws.Range("M22").Value = 22.34478233333333
$ws.Range("N22").Value = 67.034347
$ws.Range("O22").Value = 0.03165884810812076
$ws.Range("P22").Value = 0.03165884810812076
$ws.Range("Q22").Value = 605.0370301213151
$ws.Range("R22").Value = 5445.333271091836
$ws.Range("S22").Value = 0.0004370671247330308
$ws.Range("T22").Value = 0.0004370671247330308
$ws.Range("E23").Value = 3.0
$ws.Range("G23").Value = 27.07732933333333
$ws.Range("H23").Value = 81.231988
$ws.Range("I23").Value = 0.01380552833888228
$ws.Range("J23").Value = 0.01380552833888228
$ws.Range("K23").Value = 3.0
$ws.Range("M23").Value = 92.44713066666667
$ws.Range("N23").Value = 277.341392
$ws.Range("O23").Value = 0.1309822411400946
$ws.Range("P23").Value = 0.1309822411400946
$ws.Range("Q23").Value = 2503.221402983033
$ws.Range("R23").Value = 22528.9926268473
$ws.Range("S23").Value = 0.001808279041949889
$ws.Range("T23").Value = 0.001808279041949888
$ws.Range("E24").Value = 3.0
$ws.Range("G24").Value = 27.07732933333333
$ws.Range("H24").Value = 81.231988
$ws.Range("I24").Value = 0.01380552833888228
$ws.Range("J24").Value = 0.01380552833888228
$ws.Range("K24").Value = 3.0
$ws.Range("M24").Value = 243.96462
$ws.Range("N24").Value = 731.89386
$ws.Range("O24").Value = 0.3456573768818275
$ws.Range("P24").Value = 0.3456573768818275
$ws.Range("Q24").Value = 6605.91036142152
$ws.Range("R24").Value = 59453.19325279368
$ws.Range("S24").Value = 0.004771982712085781
$ws.Range("T24").Value = 0.004771982712085781
$ws.Range("E25").Value = 3.0
$ws.Range("G25").Value = 27.07732933333333
$ws.Range("H25").Value = 81.231988
$ws.Range("I25").Value = 0.01380552833888228
$ws.Range("J25").Value = 0.01380552833888228
$ws.Range("K25").Value = 3.0
$ws.Range("M25").Value = 281.5837096666667
$ws.Range("N25").Value = 844.751129
$ws.Range("O25").Value = 0.3989573834764815
$ws.Range("P25").Value = 0.3989573834764815
$ws.Range("Q25").Value = 7624.53484154605
$ws.Range("R25").Value = 68620.81357391445
$ws.Range("S25").Value = 0.00550781746359089
$ws.Range("T25").Value = 0.00550781746359089
$ws.Range("E26").Value = 3.0
$ws.Range("G26").Value = 27.07732933333333
$ws.Range("H26").Value = 81.231988
$ws.Range("I26").Value = 0.01380552833888228
$ws.Range("J26").Value = 0.01380552833888228
$ws.Range("K26").Value = 3.0
$ws.Range("M26").Value = 65.45872566666667
$ws.Range("N26").Value = 196.376177
$ws.Range("O26").Value = 0.09274415039347572
$ws.Range("P26").Value = 0.09274415039347571
$ws.Range("Q26").Value = 1772.447472616653
$ws.Range("R26").Value = 15952.02725354988
$ws.Range("S26").Value = 0.001280381996522689
$ws.Range("T26").Value = 0.001280381996522689

Write-Output "Done"